# Auto update Excel log
# Appends newly-logged sensor events to the Proximity, mmWave and Camera
# sheets of the SeniorConnect master log.
#
# Note: date-shaped strings (e.g. "2026-02-01") are written with a leading
# apostrophe so Excel stores them as literal text (matching the existing
# rows) instead of auto-converting them to a date serial number.

$wb = $excel.ActiveWorkbook

# Positional parameters only -- this interpreter does not reliably bind
# PowerShell named (-Param value) arguments on user-defined functions.
function Add-LogRow {
    param($Sheet, $Row, $Date, $Timestamp, $Hour, $Location, $Value, $Status)

    $Sheet.Cells.Item($Row, 1).Value = "'" + $Date
    $Sheet.Cells.Item($Row, 2).Value = $Timestamp
    $Sheet.Cells.Item($Row, 3).Value = $Hour
    $Sheet.Cells.Item($Row, 4).Value = $Location
    $Sheet.Cells.Item($Row, 5).Value = $Value
    $Sheet.Cells.Item($Row, 6).Value = $Status
}

# --- Proximity sheet: two new EXIT events on the Living Room Main Door ---
$wsProximity = $wb.Worksheets.Item("Proximity")

Add-LogRow $wsProximity 32 "2026-02-01" "14:41:24" "14:00" "Living Room Main Door" "EXIT" "User EXITED Living Room Main Door"
Add-LogRow $wsProximity 33 "2026-02-01" "14:41:25" "14:00" "Living Room Main Door" "EXIT" "User EXITED Living Room Main Door"

# --- mmWave sheet: seven new PRESENCE_DETECTED events in the Living Room ---
$wsMmWave = $wb.Worksheets.Item("mmWave")

Add-LogRow $wsMmWave 16 "2026-02-01" "14:41:13" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $wsMmWave 17 "2026-02-01" "14:41:24" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $wsMmWave 18 "2026-02-01" "14:41:29" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $wsMmWave 19 "2026-02-01" "14:41:34" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $wsMmWave 20 "2026-02-01" "14:41:45" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $wsMmWave 21 "2026-02-01" "14:41:55" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $wsMmWave 22 "2026-02-01" "14:42:06" "14:00" "Living Room" "PRESENCE_DETECTED" "Active"

# --- Camera sheet: image captured / received events on the Main Door ---
$wsCamera = $wb.Worksheets.Item("Camera")

Add-LogRow $wsCamera 20 "2026-02-01" "14:41:24" "14:00" "Living Room Main Door" "Image Captured" "Active"
Add-LogRow $wsCamera 21 "2026-02-01" "14:41:25" "14:00" "Living Room Main Door" "Image Received" "Active"
